$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has headers in B1:G1 (TB, d2S, K, IP, Win, sum) with
# data in row 2. We add a new "Save" column in H, extending the header row
# and adding a data value of 0 for the single data row.

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell H1 so it reuses the same bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header label and the new data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
